$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.917.92"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.903.56"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.84"
$ws.Range("E5").Value = "  -3.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.68"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "2.901.60"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.05"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.18"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "3.383.91"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "61.879.08"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "2.895.29"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.96"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("E21").Value = "  -4.51%  "
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.86"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.93"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.01"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.18"
$ws.Range("E26").Value = "  -9.26%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000111"
$ws.Range("E29").Value = "  +7.58%  "
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("E32").Value = "  -6.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.58"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.958"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.99"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.15"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("D45").Value = "2.703.24"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0335"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.71"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "345.58"
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.57"
$ws.Range("E51").Value = "  -5.32%  "
